$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List1")

# --- Row 3 headers: add two new measurement columns, rename B3/C3 ---
$ws.Range("B3").Value = "Čištění [s]"
$ws.Range("C3").Value = "Výměna barev [s]"
$ws.Range("D3").Value = "Vyschnutí aplikátoru  [s]"
$ws.Range("E3").Value = "Výměna celkem  [s]"

# --- Row 4: Primer values, now numeric with a drying-time + total formula ---
$ws.Range("B4").Value = 88
$ws.Range("C4").Value = 254
$ws.Range("D4").Value = 12
$ws.Range("E4").Formula = "=C4+D4"

# --- Row 7: Base values, now numeric with a drying-time + total formula ---
$ws.Range("B7").Value = 40
$ws.Range("C7").Value = 100
$ws.Range("D7").Value = 12
$ws.Range("E7").Formula = "=C7+D7"

# --- Row 10: Clear values, now numeric with a drying-time + total formula ---
$ws.Range("B10").Value = 88
$ws.Range("C10").Value = 254
$ws.Range("D10").Value = 12
$ws.Range("E10").Formula = "=C10+D10"

# --- Column widths (B/C/D) ---
$ws.Columns.Item(2).ColumnWidth = 12.5
$ws.Columns.Item(3).ColumnWidth = 15
$ws.Columns.Item(4).ColumnWidth = 20.833333333333332

# --- Selection moves to A13 ---
$ws.Range("A13").Select()
